# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has now
# been handed off (was "Handed back: in sync with en-US" / stale handback
# file references), for both the Overview sheet and the per-locale
# (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d6948608a1b55c438c55fa4510adc2db3a3e5c6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/555f3f8fb67b366ea6ec5029b37dd5ba0996f43f/e2e/b.md."

# --- Overview sheet: row for b.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 22:42:17"

# --- zh-cn sheet: row for b.md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-04 22:42:13"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row for b.md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-04 22:42:17"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
